# ---------------------------------------------------------------------------
# Harvard case classification was added to the underlying query-level data.
# This reclassifies which queries count as failures, which changes every
# '<provider>_old' baseline column (average/variance/std Dev) and also swaps
# which of the two rightmost summary columns is the live "average_doctor" vs
# the "average_doctor_old" baseline (the old average_doctor values become the
# new average_doctor_old values, and average_doctor is recomputed).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "average_doctor" / "average_doctor_old" column headers (BP <-> BQ)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4 (stats_for_precision): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E4").Value = 0.441
$ws.Range("F4").Value = 0.067
$ws.Range("G4").Value = 0.258
$ws.Range("N4").Value = 0.445
$ws.Range("O4").Value = 0.068
$ws.Range("P4").Value = 0.261
$ws.Range("Q4").Value = 0.018
$ws.Range("R4").Value = 0.013
$ws.Range("S4").Value = 0.116
$ws.Range("W4").Value = 0.291
$ws.Range("X4").Value = 0.112
$ws.Range("Y4").Value = 0.334
$ws.Range("AI4").Value = 0.323
$ws.Range("AJ4").Value = 0.08699999999999999
$ws.Range("AK4").Value = 0.296
$ws.Range("AU4").Value = 0.194
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.171
$ws.Range("BA4").Value = 2.013
$ws.Range("BB4").Value = 0.153
$ws.Range("BC4").Value = 0.391
$ws.Range("BG4").Value = 0.73
$ws.Range("BH4").Value = 0.138
$ws.Range("BI4").Value = 0.372
$ws.Range("BM4").Value = 0.723
$ws.Range("BN4").Value = 0.075
$ws.Range("BO4").Value = 0.274
$ws.Range("BP4").Value = 0.671
$ws.Range("BQ4").Value = 0.705

# Row 5 (stats_for_recall): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E5").Value = 0.5649999999999999
$ws.Range("F5").Value = 0.074
$ws.Range("G5").Value = 0.272
$ws.Range("N5").Value = 0.731
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.282
$ws.Range("Q5").Value = 0.008999999999999999
$ws.Range("R5").Value = 0.002
$ws.Range("S5").Value = 0.05
$ws.Range("W5").Value = 0.272
$ws.Range("X5").Value = 0.103
$ws.Range("Y5").Value = 0.321
$ws.Range("AI5").Value = 0.346
$ws.Range("AJ5").Value = 0.097
$ws.Range("AK5").Value = 0.312
$ws.Range("AU5").Value = 0.373
$ws.Range("AV5").Value = 0.096
$ws.Range("AW5").Value = 0.31
$ws.Range("BA5").Value = 1.344
$ws.Range("BB5").Value = 0.078
$ws.Range("BC5").Value = 0.278
$ws.Range("BG5").Value = 0.398
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.227
$ws.Range("BM5").Value = 0.553
$ws.Range("BN5").Value = 0.061
$ws.Range("BO5").Value = 0.248
$ws.Range("BP5").Value = 0.448
$ws.Range("BQ5").Value = 0.455

# Row 6 (stats_for_f1-score): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E6").Value = 0.495
$ws.Range("N6").Value = 0.553
$ws.Range("Q6").Value = 0.012
$ws.Range("W6").Value = 0.281
$ws.Range("AI6").Value = 0.334
$ws.Range("AU6").Value = 0.255
$ws.Range("BA6").Value = 1.604
$ws.Range("BG6").Value = 0.515
$ws.Range("BM6").Value = 0.627
$ws.Range("BP6").Value = 0.535
$ws.Range("BQ6").Value = 0.55

# Row 7 (stats_for_f2-score): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E7").Value = 0.535
$ws.Range("N7").Value = 0.648
$ws.Range("Q7").Value = 0.01
$ws.Range("W7").Value = 0.276
$ws.Range("AI7").Value = 0.341
$ws.Range("AU7").Value = 0.315
$ws.Range("BA7").Value = 1.436
$ws.Range("BG7").Value = 0.438
$ws.Range("BM7").Value = 0.58
$ws.Range("BP7").Value = 0.479
$ws.Range("BQ7").Value = 0.489

# Row 8 (stats_for_NDCG): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E8").Value = 0.634
$ws.Range("F8").Value = 0.096
$ws.Range("G8").Value = 0.31
$ws.Range("N8").Value = 0.778
$ws.Range("O8").Value = 0.064
$ws.Range("P8").Value = 0.253
$ws.Range("Q8").Value = 0.01
$ws.Range("R8").Value = 0.006
$ws.Range("W8").Value = 0.308
$ws.Range("X8").Value = 0.122
$ws.Range("Y8").Value = 0.349
$ws.Range("AI8").Value = 0.37
$ws.Range("AJ8").Value = 0.13
$ws.Range("AK8").Value = 0.36
$ws.Range("AU8").Value = 0.316
$ws.Range("AV8").Value = 0.08500000000000001
$ws.Range("AW8").Value = 0.291
$ws.Range("BA8").Value = 1.739
$ws.Range("BB8").Value = 0.12
$ws.Range("BC8").Value = 0.346
$ws.Range("BG8").Value = 0.5580000000000001
$ws.Range("BH8").Value = 0.105
$ws.Range("BI8").Value = 0.324
$ws.Range("BM8").Value = 0.694
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.249
$ws.Range("BP8").Value = 0.58
$ws.Range("BQ8").Value = 0.599

# Row 9 (stats_for_M1): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E9").Value = 0.571
$ws.Range("F9").Value = 0.245
$ws.Range("G9").Value = 0.495
$ws.Range("N9").Value = 0.6899999999999999
$ws.Range("O9").Value = 0.214
$ws.Range("P9").Value = 0.462
$ws.Range("W9").Value = 0.214
$ws.Range("X9").Value = 0.168
$ws.Range("Y9").Value = 0.41
$ws.Range("AI9").Value = 0.286
$ws.Range("AJ9").Value = 0.204
$ws.Range("AK9").Value = 0.452
$ws.Range("BA9").Value = 1.666
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.583
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.643
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.555
$ws.Range("BQ9").Value = 0.5679999999999999

# Row 10 (stats_for_M3): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E10").Value = 0.714
$ws.Range("F10").Value = 0.204
$ws.Range("G10").Value = 0.452
$ws.Range("N10").Value = 0.893
$ws.Range("O10").Value = 0.096
$ws.Range("P10").Value = 0.309
$ws.Range("W10").Value = 0.381
$ws.Range("X10").Value = 0.236
$ws.Range("Y10").Value = 0.486
$ws.Range("AI10").Value = 0.405
$ws.Range("AJ10").Value = 0.241
$ws.Range("AK10").Value = 0.491
$ws.Range("AU10").Value = 0.31
$ws.Range("AV10").Value = 0.214
$ws.Range("AW10").Value = 0.462
$ws.Range("BA10").Value = 2.083
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.643
$ws.Range("BH10").Value = 0.23
$ws.Range("BI10").Value = 0.479
$ws.Range("BM10").Value = 0.857
$ws.Range("BN10").Value = 0.122
$ws.Range("BO10").Value = 0.35
$ws.Range("BP10").Value = 0.694
$ws.Range("BQ10").Value = 0.721

# Row 11 (stats_for_M5): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E11").Value = 0.75
$ws.Range("F11").Value = 0.188
$ws.Range("G11").Value = 0.433
$ws.Range("N11").Value = 0.905
$ws.Range("O11").Value = 0.08599999999999999
$ws.Range("P11").Value = 0.294
$ws.Range("W11").Value = 0.381
$ws.Range("X11").Value = 0.236
$ws.Range("Y11").Value = 0.486
$ws.Range("AI11").Value = 0.44
$ws.Range("AJ11").Value = 0.246
$ws.Range("AK11").Value = 0.496
$ws.Range("AU11").Value = 0.44
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.083
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.643
$ws.Range("BH11").Value = 0.23
$ws.Range("BI11").Value = 0.479
$ws.Range("BM11").Value = 0.857
$ws.Range("BN11").Value = 0.122
$ws.Range("BO11").Value = 0.35
$ws.Range("BP11").Value = 0.694
$ws.Range("BQ11").Value = 0.724

# Row 12 (stats_for_position): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.785
$ws.Range("G12").Value = 0.886
$ws.Range("N12").Value = 1.39
$ws.Range("O12").Value = 0.705
$ws.Range("P12").Value = 0.84
$ws.Range("W12").Value = 1.594
$ws.Range("X12").Value = 0.554
$ws.Range("Y12").Value = 0.744
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.769
$ws.Range("AV12").Value = 2.844
$ws.Range("AW12").Value = 1.686
$ws.Range("BA12").Value = 3.777
$ws.Range("BB12").Value = 0.431
$ws.Range("BC12").Value = 0.656
$ws.Range("BG12").Value = 1.111
$ws.Range("BH12").Value = 0.136
$ws.Range("BI12").Value = 0.369
$ws.Range("BM12").Value = 1.319
$ws.Range("BN12").Value = 0.356
$ws.Range("BO12").Value = 0.597
$ws.Range("BP12").Value = 1.259
$ws.Range("BQ12").Value = 1.282

# Row 13 (stats_for_length (x of gs)): updated "_old" columns + average_doctor / average_doctor_old
$ws.Range("E13").Value = 1.561
$ws.Range("F13").Value = 0.667
$ws.Range("G13").Value = 0.8169999999999999
$ws.Range("N13").Value = 2.022
$ws.Range("O13").Value = 0.994
$ws.Range("P13").Value = 0.997
$ws.Range("W13").Value = 1.026
$ws.Range("X13").Value = 0.188
$ws.Range("Y13").Value = 0.434
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.278
$ws.Range("AV13").Value = 0.997
$ws.Range("AW13").Value = 0.999
$ws.Range("BA13").Value = 2.35
$ws.Range("BB13").Value = 0.304
$ws.Range("BC13").Value = 0.551
$ws.Range("BG13").Value = 0.59
$ws.Range("BH13").Value = 0.075
$ws.Range("BI13").Value = 0.274
$ws.Range("BM13").Value = 0.889
$ws.Range("BN13").Value = 0.285
$ws.Range("BO13").Value = 0.533
$ws.Range("BP13").Value = 0.783
$ws.Range("BQ13").Value = 0.728
